$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the expected-result assertion for the Facebook login scenario (D2)
$ws.Range("D2").Value = "Facebook"

# Move the active selection to reflect the asserted cell
$ws.Range("D2").Select()
